# Apply cryptos list update (price/volume refresh + two row-pair swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "593.17") are stored verbatim instead of being coerced to a
# binary double (which would introduce float noise / sci-notation).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.696.26"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").Value = "2.734.11"
$ws.Range("E3").Value = "  +4.64%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "593.17"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "152.92"
$ws.Range("E6").Value = "  +6.75%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +2.14%  "
$ws.Range("D9").Value = "2.765.10"
$ws.Range("E9").Value = "  +5.43%  "
$ws.Range("D11").Value = "0.113"
$ws.Range("E11").Value = "  +7.43%  "
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("D14").Value = "3.229.43"
$ws.Range("E14").Value = "  +5.12%  "
$ws.Range("D15").Value = "26.59"
$ws.Range("E15").Value = "  +6.88%  "
$ws.Range("D16").Value = "63.622.90"
$ws.Range("E16").Value = "  +5.52%  "
$ws.Range("E17").Value = "  +8.75%  "
$ws.Range("D18").Value = "2.763.25"
$ws.Range("E18").Value = "  +5.62%  "
$ws.Range("D19").Value = "12.05"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "4.91"
$ws.Range("E20").Value = "  +5.74%  "
$ws.Range("D21").Value = "365.66"
$ws.Range("E21").Value = "  +5.40%  "
$ws.Range("D22").Value = "7.02"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.537"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "65.94"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +5.16%  "
$ws.Range("D27").Value = "8.67"
$ws.Range("E27").Value = "  +8.19%  "
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("E29").Value = "  +14.01%  "
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +4.70%  "
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("D32").Value = "172.46"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +18.48%  "
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "20.64"
$ws.Range("E35").Value = "  +5.84%  "
$ws.Range("E36").Value = "  +12.63%  "
$ws.Range("D37").Value = "1.43"
$ws.Range("E37").Value = "  +9.23%  "
$ws.Range("D38").Value = "1.79"
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("D39").Value = "1.01"
$ws.Range("E39").Value = "  +19.22%  "
$ws.Range("D40").Value = "348.15"
$ws.Range("E40").Value = "  +9.06%  "
$ws.Range("E41").Value = "  +7.60%  "
$ws.Range("D42").Value = "38.95"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "22.15"
$ws.Range("E43").Value = "  +10.72%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "5.60"
$ws.Range("E44").Value = "  +11.70%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "143.34"
$ws.Range("E45").Value = "  +5.79%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "22.23"
$ws.Range("E46").Value = "  +10.90%  "
$ws.Range("D47").Value = "0.0592"
$ws.Range("E47").Value = "  +7.28%  "
$ws.Range("D48").Value = "0.647"
$ws.Range("E48").Value = "  +6.20%  "
$ws.Range("E49").Value = "  +7.40%  "
$ws.Range("E50").Value = "  +2.53%  "
$ws.Range("D51").Value = "2.172.55"
$ws.Range("E51").Value = "  +7.34%  "

# Restore the default (unstyled) cell style now that the literal text
# is safely stored — matches the workbook's original unstyled D column.
$priceRange.Style = "Normal"

